# Auto-generated script applying the scheduled-runner price update
# Each block sets the final (post-edit) values for the affected Leve rows,
# matching currentAveragePrice/NQ/HQ, LevePrice NQ/HQ and LeveProfit NQ/HQ columns.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1358.1818
$ws.Range("I40").Value = 868
$ws.Range("K40").Value = 868
$ws.Range("M40").Value = -693
# Row 80
$ws.Range("H80").Value = 1620
$ws.Range("I80").Value = 2700
$ws.Range("J80").Value = 900
$ws.Range("K80").Value = 8100
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = -7102
$ws.Range("N80").Value = -4696
# Row 83
$ws.Range("H83").Value = 1620
$ws.Range("I83").Value = 2700
$ws.Range("J83").Value = 900
$ws.Range("K83").Value = 24300
$ws.Range("L83").Value = 8100
$ws.Range("M83").Value = -19308
$ws.Range("N83").Value = -18084
# Row 99
$ws.Range("H99").Value = 1728.6364
$ws.Range("I99").Value = 378.66666
$ws.Range("J99").Value = 3348.6
$ws.Range("K99").Value = 1135.99998
$ws.Range("L99").Value = 10045.8
$ws.Range("M99").Value = 362.0000199999999
$ws.Range("N99").Value = -13041.8
# Row 137
$ws.Range("H137").Value = 1571.0938
$ws.Range("I137").Value = 1333.1364
$ws.Range("J137").Value = 2094.6
$ws.Range("K137").Value = 3999.4092
$ws.Range("L137").Value = 6283.799999999999
$ws.Range("M137").Value = -1449.4092
$ws.Range("N137").Value = -11383.8
# Row 138
$ws.Range("H138").Value = 574016.6
$ws.Range("I138").Value = 1884.375
$ws.Range("K138").Value = 5653.125
$ws.Range("M138").Value = -513.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3218.75
$ws.Range("I45").Value = 3107.1428
$ws.Range("K45").Value = 3107.1428
$ws.Range("M45").Value = -2730.1428
# Row 61
$ws.Range("H61").Value = 2211.4
$ws.Range("I61").Value = 2100
$ws.Range("K61").Value = 2100
$ws.Range("M61").Value = -1888
# Row 74
$ws.Range("H74").Value = 1466
$ws.Range("I74").Value = 1608
$ws.Range("J74").Value = 614
$ws.Range("K74").Value = 1608
$ws.Range("L74").Value = 614
$ws.Range("M74").Value = -734
$ws.Range("N74").Value = -2362
# Row 77
$ws.Range("H77").Value = 1466
$ws.Range("I77").Value = 1608
$ws.Range("J77").Value = 614
$ws.Range("K77").Value = 8040
$ws.Range("L77").Value = 3070
$ws.Range("M77").Value = -3672
$ws.Range("N77").Value = -11806
# Row 132
$ws.Range("H132").Value = 2912.2856
$ws.Range("I132").Value = 2796.2144
$ws.Range("J132").Value = 3144.4285
$ws.Range("K132").Value = 8388.643199999999
$ws.Range("L132").Value = 9433.2855
$ws.Range("M132").Value = -5858.643199999999
$ws.Range("N132").Value = -14493.2855
# Row 136
$ws.Range("H136").Value = 2211.4
$ws.Range("I136").Value = 2100
$ws.Range("K136").Value = 6300
$ws.Range("M136").Value = -3750

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2694.25
$ws.Range("I86").Value = 2846.6
$ws.Range("J86").Value = 2367.7856
$ws.Range("K86").Value = 2846.6
$ws.Range("L86").Value = 2367.7856
$ws.Range("M86").Value = -1723.6
$ws.Range("N86").Value = -4613.7856
# Row 89
$ws.Range("H89").Value = 2694.25
$ws.Range("I89").Value = 2846.6
$ws.Range("J89").Value = 2367.7856
$ws.Range("K89").Value = 14233
$ws.Range("L89").Value = 11838.928
$ws.Range("M89").Value = -8617
$ws.Range("N89").Value = -23070.928
# Row 105
$ws.Range("H105").Value = 200002000
$ws.Range("I105").Value = 200002000
$ws.Range("K105").Value = 200002000
$ws.Range("M105").Value = -200000253
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
# Row 132
$ws.Range("H132").Value = 49165.332
$ws.Range("J132").Value = 49165.332
$ws.Range("L132").Value = 49165.332
$ws.Range("N132").Value = -59285.332

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 48499.25
$ws.Range("J20").Value = 48499.25
$ws.Range("L20").Value = 48499.25
$ws.Range("N20").Value = -48971.25
# Row 30
$ws.Range("H30").Value = 48499.25
$ws.Range("J30").Value = 48499.25
$ws.Range("L30").Value = 48499.25
$ws.Range("N30").Value = -48681.25
# Row 31
$ws.Range("H31").Value = 1597.0588
$ws.Range("I31").Value = 1275.1428
$ws.Range("K31").Value = 1275.1428
$ws.Range("M31").Value = -980.1428000000001
# Row 34
$ws.Range("H34").Value = 1597.0588
$ws.Range("I34").Value = 1275.1428
$ws.Range("K34").Value = 1275.1428
$ws.Range("M34").Value = -1073.1428
# Row 70
$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630
# Row 73
$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184
# Row 107
$ws.Range("H107").Value = 823.9091
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 823.9091
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 823.9091
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = -4663.9091
# Row 128
$ws.Range("H128").Value = 48499.25
$ws.Range("J128").Value = 48499.25
$ws.Range("L128").Value = 48499.25
$ws.Range("N128").Value = -58459.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 2207.1428
# Row 80
$ws.Range("H80").Value = 5252.5454
$ws.Range("J80").Value = 5252.5454
$ws.Range("L80").Value = 15757.6362
$ws.Range("N80").Value = -17629.6362
# Row 83
$ws.Range("H83").Value = 5252.5454
$ws.Range("J83").Value = 5252.5454
$ws.Range("L83").Value = 47272.9086
$ws.Range("N83").Value = -56632.9086

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2433.0715
$ws.Range("I126").Value = 1337.2
$ws.Range("K126").Value = 4011.6
$ws.Range("M126").Value = -1541.6
# Row 132
$ws.Range("H132").Value = 3475.2778
$ws.Range("I132").Value = 2682.7144
$ws.Range("J132").Value = 6249.25
$ws.Range("K132").Value = 8048.1432
$ws.Range("L132").Value = 18747.75
$ws.Range("M132").Value = -5518.1432
$ws.Range("N132").Value = -23807.75
# Row 140
$ws.Range("H140").Value = 42000
$ws.Range("J140").Value = 42000
$ws.Range("L140").Value = 42000
$ws.Range("N140").Value = -52360
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = ""

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3078.3333
$ws.Range("I7").Value = 2916.25
$ws.Range("K7").Value = 2916.25
$ws.Range("M7").Value = -2804.25
# Row 68
$ws.Range("H68").Value = 1474.421
$ws.Range("I68").Value = 1310.2727
$ws.Range("K68").Value = 1310.2727
$ws.Range("M68").Value = -561.2727
# Row 71
$ws.Range("H71").Value = 1474.421
$ws.Range("I71").Value = 1310.2727
$ws.Range("K71").Value = 6551.363499999999
$ws.Range("M71").Value = -2807.363499999999
# Row 100
$ws.Range("H100").Value = 2201.5
$ws.Range("I100").Value = 2002
$ws.Range("K100").Value = 2002
$ws.Range("M100").Value = -1461
# Row 122
$ws.Range("H122").Value = 35422692
$ws.Range("I122").Value = 56670508
$ws.Range("J122").Value = 9666.666999999999
$ws.Range("K122").Value = 170011524
$ws.Range("L122").Value = 29000.001
$ws.Range("M122").Value = -170009074
$ws.Range("N122").Value = -33900.001
# Row 126
$ws.Range("H126").Value = 3078.3333
$ws.Range("I126").Value = 2916.25
$ws.Range("K126").Value = 8748.75
$ws.Range("M126").Value = -6278.75
# Row 132
$ws.Range("H132").Value = 34287.97
$ws.Range("I132").Value = 1564.7727
$ws.Range("K132").Value = 4694.3181
$ws.Range("M132").Value = -2164.3181
# Row 136
$ws.Range("H136").Value = 7983.5625
$ws.Range("I136").Value = 15596.571
$ws.Range("J136").Value = 2062.3333
$ws.Range("K136").Value = 46789.713
$ws.Range("L136").Value = 6186.999899999999
$ws.Range("M136").Value = -44239.713
$ws.Range("N136").Value = -11286.9999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 28670000
$ws.Range("J5").Value = 28670000
$ws.Range("L5").Value = 28670000
$ws.Range("N5").Value = -28670224
# Row 56
$ws.Range("H56").Value = 30000
$ws.Range("J56").Value = 30000
$ws.Range("L56").Value = 30000
$ws.Range("N56").Value = -31428
# Row 62
$ws.Range("H62").Value = 35722590
$ws.Range("I62").Value = 55561424
$ws.Range("K62").Value = 55561424
$ws.Range("M62").Value = -55560800
# Row 65
$ws.Range("H65").Value = 35722590
$ws.Range("I65").Value = 55561424
$ws.Range("K65").Value = 277807120
$ws.Range("M65").Value = -277804000
# Row 132
$ws.Range("H132").Value = 4651.647
$ws.Range("I132").Value = 5171.136
$ws.Range("K132").Value = 15513.408
$ws.Range("M132").Value = -12983.408
